$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.692291333333333
$ws.Range("H2").Value = 8.076874
$ws.Range("I2").Value = 0.1057975873398775
$ws.Range("J2").Value = 0.1132142695545834
$ws.Range("M2").Value = 63.46725166666666
$ws.Range("N2").Value = 190.401755
$ws.Range("O2").Value = 0.2354497988808272
$ws.Range("P2").Value = 0.2397164477183668
$ws.Range("Q2").Value = 170.8723316126522
$ws.Range("R2").Value = 1537.85098451387
$ws.Range("S2").Value = 0.02491002066125091
$ws.Range("T2").Value = 0.02713932252865437
$ws.Range("G3").Value = 2.692291333333333
$ws.Range("H3").Value = 8.076874
$ws.Range("I3").Value = 0.1057975873398775
$ws.Range("J3").Value = 0.1132142695545834
$ws.Range("M3").Value = 47.980825
$ws.Range("N3").Value = 143.942475
$ws.Range("O3").Value = 0.1779985000094065
$ws.Range("P3").Value = 0.1812240584798697
$ws.Range("Q3").Value = 129.1783593136833
$ws.Range("R3").Value = 1162.60523382315
$ws.Range("S3").Value = 0.01883181185111236
$ws.Range("T3").Value = 0.02051714940651555
$ws.Range("G4").Value = 2.692291333333333
$ws.Range("H4").Value = 8.076874
$ws.Range("I4").Value = 0.1057975873398775
$ws.Range("J4").Value = 0.1132142695545834
$ws.Range("M4").Value = 64.53809233333334
$ws.Range("N4").Value = 193.614277
$ws.Range("O4").Value = 0.2394223865221556
$ws.Range("P4").Value = 0.243761023683841
$ws.Range("Q4").Value = 173.7553466588998
$ws.Range("R4").Value = 1563.798119930098
$ws.Range("S4").Value = 0.02533031084919967
$ws.Range("T4").Value = 0.02759722624224356
$ws.Range("G5").Value = 2.692291333333333
$ws.Range("H5").Value = 8.076874
$ws.Range("I5").Value = 0.1057975873398775
$ws.Range("J5").Value = 0.1132142695545834
$ws.Range("M5").Value = 14.3933435
$ws.Range("N5").Value = 28.786687
$ws.Range("O5").Value = 0.0533961963580272
$ws.Range("P5").Value = 0.03624253541791403
$ws.Range("Q5").Value = 38.75107396273967
$ws.Range("R5").Value = 232.506443776438
$ws.Range("S5").Value = 0.00564918874780563
$ws.Range("T5").Value = 0.004103172174145255
$ws.Range("G6").Value = 2.692291333333333
$ws.Range("H6").Value = 8.076874
$ws.Range("I6").Value = 0.1057975873398775
$ws.Range("J6").Value = 0.1132142695545834
$ws.Range("M6").Value = 79.17795566666666
$ws.Range("N6").Value = 237.533867
$ws.Range("O6").Value = 0.2937331182295834
$ws.Range("P6").Value = 0.2990559347000084
$ws.Range("Q6").Value = 213.1701238324175
$ws.Range("R6").Value = 1918.531114491758
$ws.Range("S6").Value = 0.03107625523050891
$ws.Range("T6").Value = 0.03385739920302464
$ws.Range("G7").Value = 4.524801333333333
$ws.Range("H7").Value = 13.574404
$ws.Range("I7").Value = 0.1778087899819636
$ws.Range("J7").Value = 0.1902736421911268
$ws.Range("M7").Value = 63.46725166666666
$ws.Range("N7").Value = 190.401755
$ws.Range("O7").Value = 0.2354497988808272
$ws.Range("P7").Value = 0.2397164477183668
$ws.Range("Q7").Value = 287.1767049643356
$ws.Range("R7").Value = 2584.59034467902
$ws.Range("S7").Value = 0.04186504384049657
$ws.Range("T7").Value = 0.04561172160049248
$ws.Range("G8").Value = 4.524801333333333
$ws.Range("H8").Value = 13.574404
$ws.Range("I8").Value = 0.1778087899819636
$ws.Range("J8").Value = 0.1902736421911268
$ws.Range("M8").Value = 47.980825
$ws.Range("N8").Value = 143.942475
$ws.Range("O8").Value = 0.1779985000094065
$ws.Range("P8").Value = 0.1812240584798697
$ws.Range("Q8").Value = 217.1037009344334
$ws.Range("R8").Value = 1953.9333084099
$ws.Range("S8").Value = 0.03164969790527709
$ws.Range("T8").Value = 0.03448216165962257
$ws.Range("G9").Value = 4.524801333333333
$ws.Range("H9").Value = 13.574404
$ws.Range("I9").Value = 0.1778087899819636
$ws.Range("J9").Value = 0.1902736421911268
$ws.Range("M9").Value = 64.53809233333334
$ws.Range("N9").Value = 193.614277
$ws.Range("O9").Value = 0.2394223865221556
$ws.Range("P9").Value = 0.243761023683841
$ws.Range("Q9").Value = 292.0220462406565
$ws.Range("R9").Value = 2628.198416165908
$ws.Range("S9").Value = 0.04257140484209848
$ws.Range("T9").Value = 0.04638129780056196
$ws.Range("G10").Value = 4.524801333333333
$ws.Range("H10").Value = 13.574404
$ws.Range("I10").Value = 0.1778087899819636
$ws.Range("J10").Value = 0.1902736421911268
$ws.Range("M10").Value = 14.3933435
$ws.Range("N10").Value = 28.786687
$ws.Range("O10").Value = 0.0533961963580272
$ws.Range("P10").Value = 0.03624253541791403
$ws.Range("Q10").Value = 65.12701985992467
$ws.Range("R10").Value = 390.762119159548
$ws.Range("S10").Value = 0.009494313064060147
$ws.Range("T10").Value = 0.006895999216207416
$ws.Range("G11").Value = 4.524801333333333
$ws.Range("H11").Value = 13.574404
$ws.Range("I11").Value = 0.1778087899819636
$ws.Range("J11").Value = 0.1902736421911268
$ws.Range("M11").Value = 79.17795566666666
$ws.Range("N11").Value = 237.533867
$ws.Range("O11").Value = 0.2937331182295834
$ws.Range("P11").Value = 0.2990559347000084
$ws.Range("Q11").Value = 358.2645193711409
$ws.Range("R11").Value = 3224.380674340268
$ws.Range("S11").Value = 0.05222833033003127
$ws.Range("T11").Value = 0.05690246191424238
$ws.Range("G12").Value = 4.544410333333333
$ws.Range("H12").Value = 13.633231
$ws.Range("I12").Value = 0.1785793547661169
$ws.Range("J12").Value = 0.1910982255429393
$ws.Range("M12").Value = 63.46725166666666
$ws.Range("N12").Value = 190.401755
$ws.Range("O12").Value = 0.2354497988808272
$ws.Range("P12").Value = 0.2397164477183668
$ws.Range("Q12").Value = 288.4212343022672
$ws.Range("R12").Value = 2595.791108720405
$ws.Range("S12").Value = 0.04204647316395011
$ws.Range("T12").Value = 0.04580938779243669
$ws.Range("G13").Value = 4.544410333333333
$ws.Range("H13").Value = 13.633231
$ws.Range("I13").Value = 0.1785793547661169
$ws.Range("J13").Value = 0.1910982255429393
$ws.Range("M13").Value = 47.980825
$ws.Range("N13").Value = 143.942475
$ws.Range("O13").Value = 0.1779985000094065
$ws.Range("P13").Value = 0.1812240584798697
$ws.Range("Q13").Value = 218.0445569318584
$ws.Range("R13").Value = 1962.401012386725
$ws.Range("S13").Value = 0.03178685728101645
$ws.Range("T13").Value = 0.03463159600119297
$ws.Range("G14").Value = 4.544410333333333
$ws.Range("H14").Value = 13.633231
$ws.Range("I14").Value = 0.1785793547661169
$ws.Range("J14").Value = 0.1910982255429393
$ws.Range("M14").Value = 64.53809233333334
$ws.Range("N14").Value = 193.614277
$ws.Range("O14").Value = 0.2394223865221556
$ws.Range("P14").Value = 0.243761023683841
$ws.Range("Q14").Value = 293.2875736932208
$ws.Range("R14").Value = 2639.588163238987
$ws.Range("S14").Value = 0.04275589530169039
$ws.Range("T14").Value = 0.04658229908251243
$ws.Range("G15").Value = 4.544410333333333
$ws.Range("H15").Value = 13.633231
$ws.Range("I15").Value = 0.1785793547661169
$ws.Range("J15").Value = 0.1910982255429393
$ws.Range("M15").Value = 14.3933435
$ws.Range("N15").Value = 28.786687
$ws.Range("O15").Value = 0.0533961963580272
$ws.Range("P15").Value = 0.03624253541791403
$ws.Range("Q15").Value = 65.40925893261617
$ws.Range("R15").Value = 392.455553595697
$ws.Range("S15").Value = 0.009535458292581376
$ws.Range("T15").Value = 0.006925884207540502
$ws.Range("G16").Value = 4.544410333333333
$ws.Range("H16").Value = 13.633231
$ws.Range("I16").Value = 0.1785793547661169
$ws.Range("J16").Value = 0.1910982255429393
$ws.Range("M16").Value = 79.17795566666666
$ws.Range("N16").Value = 237.533867
$ws.Range("O16").Value = 0.2937331182295834
$ws.Range("P16").Value = 0.2990559347000084
$ws.Range("Q16").Value = 359.8171199038085
$ws.Range("R16").Value = 3238.354079134277
$ws.Range("S16").Value = 0.05245467072687852
$ws.Range("T16").Value = 0.05714905845925674
$ws.Range("G17").Value = 5.001220999999999
$ws.Range("H17").Value = 10.002442
$ws.Range("I17").Value = 0.1965304085046502
$ws.Range("J17").Value = 0.1402051294587592
$ws.Range("M17").Value = 63.46725166666666
$ws.Range("N17").Value = 190.401755
$ws.Range("O17").Value = 0.2354497988808272
$ws.Range("P17").Value = 0.2397164477183668
$ws.Range("Q17").Value = 317.4137518476182
$ws.Range("R17").Value = 1904.48251108571
$ws.Range("S17").Value = 0.0462730451563867
$ws.Range("T17").Value = 0.03360947558574749
$ws.Range("G18").Value = 5.001220999999999
$ws.Range("H18").Value = 10.002442
$ws.Range("I18").Value = 0.1965304085046502
$ws.Range("J18").Value = 0.1402051294587592
$ws.Range("M18").Value = 47.980825
$ws.Range("N18").Value = 143.942475
$ws.Range("O18").Value = 0.1779985000094065
$ws.Range("P18").Value = 0.1812240584798697
$ws.Range("Q18").Value = 239.962709587325
$ws.Range("R18").Value = 1439.77625752395
$ws.Range("S18").Value = 0.03498211792006363
$ws.Range("T18").Value = 0.02540854258021187
$ws.Range("G19").Value = 5.001220999999999
$ws.Range("H19").Value = 10.002442
$ws.Range("I19").Value = 0.1965304085046502
$ws.Range("J19").Value = 0.1402051294587592
$ws.Range("M19").Value = 64.53809233333334
$ws.Range("N19").Value = 193.614277
$ws.Range("O19").Value = 0.2394223865221556
$ws.Range("P19").Value = 0.243761023683841
$ws.Range("Q19").Value = 322.7692626774057
$ws.Range("R19").Value = 1936.615576064434
$ws.Range("S19").Value = 0.0470537794283575
$ws.Range("T19").Value = 0.03417654588259259
$ws.Range("G20").Value = 5.001220999999999
$ws.Range("H20").Value = 10.002442
$ws.Range("I20").Value = 0.1965304085046502
$ws.Range("J20").Value = 0.1402051294587592
$ws.Range("M20").Value = 14.3933435
$ws.Range("N20").Value = 28.786687
$ws.Range("O20").Value = 0.0533961963580272
$ws.Range("P20").Value = 0.03624253541791403
$ws.Range("Q20").Value = 71.98429177241348
$ws.Range("R20").Value = 287.9371670896539
$ws.Range("S20").Value = 0.0104939762828376
$ws.Range("T20").Value = 0.005081389370182301
$ws.Range("G21").Value = 5.001220999999999
$ws.Range("H21").Value = 10.002442
$ws.Range("I21").Value = 0.1965304085046502
$ws.Range("J21").Value = 0.1402051294587592
$ws.Range("M21").Value = 79.17795566666666
$ws.Range("N21").Value = 237.533867
$ws.Range("O21").Value = 0.2937331182295834
$ws.Range("P21").Value = 0.2990559347000084
$ws.Range("Q21").Value = 395.9864546172022
$ws.Range("R21").Value = 2375.918727703213
$ws.Range("S21").Value = 0.05772748971700473
$ws.Range("T21").Value = 0.04192917604002491
$ws.Range("G22").Value = 8.684844333333333
$ws.Range("H22").Value = 26.054533
$ws.Range("I22").Value = 0.341283859407392
$ws.Range("J22").Value = 0.3652087332525911
$ws.Range("M22").Value = 63.46725166666666
$ws.Range("N22").Value = 190.401755
$ws.Range("O22").Value = 0.2354497988808272
$ws.Range("P22").Value = 0.2397164477183668
$ws.Range("Q22").Value = 551.2032009894905
$ws.Range("R22").Value = 4960.828808905414
$ws.Range("S22").Value = 0.08035521605874296
$ws.Range("T22").Value = 0.08754654021103572
$ws.Range("G23").Value = 8.684844333333333
$ws.Range("H23").Value = 26.054533
$ws.Range("I23").Value = 0.341283859407392
$ws.Range("J23").Value = 0.3652087332525911
$ws.Range("M23").Value = 47.980825
$ws.Range("N23").Value = 143.942475
$ws.Range("O23").Value = 0.1779985000094065
$ws.Range("P23").Value = 0.1812240584798697
$ws.Range("Q23").Value = 416.7059961099083
$ws.Range("R23").Value = 3750.353964989175
$ws.Range("S23").Value = 0.06074801505193693
$ws.Range("T23").Value = 0.0661846088323267
$ws.Range("G24").Value = 8.684844333333333
$ws.Range("H24").Value = 26.054533
$ws.Range("I24").Value = 0.341283859407392
$ws.Range("J24").Value = 0.3652087332525911
$ws.Range("M24").Value = 64.53809233333334
$ws.Range("N24").Value = 193.614277
$ws.Range("O24").Value = 0.2394223865221556
$ws.Range("P24").Value = 0.243761023683841
$ws.Range("Q24").Value = 560.5032854852934
$ws.Range("R24").Value = 5044.529569367642
$ws.Range("S24").Value = 0.08171099610080962
$ws.Range("T24").Value = 0.08902365467593043
$ws.Range("G25").Value = 8.684844333333333
$ws.Range("H25").Value = 26.054533
$ws.Range("I25").Value = 0.341283859407392
$ws.Range("J25").Value = 0.3652087332525911
$ws.Range("M25").Value = 14.3933435
$ws.Range("N25").Value = 28.786687
$ws.Range("O25").Value = 0.0533961963580272
$ws.Range("P25").Value = 0.03624253541791403
$ws.Range("Q25").Value = 125.0039477336952
$ws.Range("R25").Value = 750.0236864021709
$ws.Range("S25").Value = 0.01822325997074245
$ws.Range("T25").Value = 0.01323609044983855
$ws.Range("G26").Value = 8.684844333333333
$ws.Range("H26").Value = 26.054533
$ws.Range("I26").Value = 0.341283859407392
$ws.Range("J26").Value = 0.3652087332525911
$ws.Range("M26").Value = 79.17795566666666
$ws.Range("N26").Value = 237.533867
$ws.Range("O26").Value = 0.2937331182295834
$ws.Range("P26").Value = 0.2990559347000084
$ws.Range("Q26").Value = 687.6482195965677
$ws.Range("R26").Value = 6188.83397636911
$ws.Range("S26").Value = 0.10024637222516
$ws.Range("T26").Value = 0.1092178390834597
